$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions data pull): prices/volumes updated
# for every existing coin, plus a new "LEO" entry inserted at rank 23 which
# pushes PEPE and every coin below it down by one row.
#
# All target cells hold plain text in the source file (e.g. "172.96",
# "0.0000178"), so each value below is written with a leading apostrophe --
# Excel's standard "force text" marker -- to stop it from silently
# re-typing numeric-looking strings as numbers. The apostrophe itself is
# not stored; Excel strips it and just keeps the cell formatted as text.
$quote = [char]39
$updates = [ordered]@{
    'D2' = '66.917.99'
    'E2' = '  +0.63%  '
    'D3' = '3.498.95'
    'E3' = '  +0.18%  '
    'E4' = '  +0.05%  '
    'D5' = '594.37'
    'D6' = '172.96'
    'E6' = '  +2.81%  '
    'E7' = '  +0.04%  '
    'E8' = '  -1.30%  '
    'E9' = '  +4.60%  '
    'D10' = '7.13'
    'E10' = '  -2.73%  '
    'E11' = '  -0.04%  '
    'D12' = '4.104.17'
    'E12' = '  +0.27%  '
    'E13' = '  +0.22%  '
    'D14' = '29.33'
    'E14' = '  +4.77%  '
    'D15' = '66.929.75'
    'E15' = '  +0.60%  '
    'D16' = '0.0000178'
    'E16' = '  +0.89%  '
    'D17' = '3.515.13'
    'E17' = '  +0.91%  '
    'D18' = '6.26'
    'E18' = '  -0.05%  '
    'E19' = '  +2.73%  '
    'D20' = '391.43'
    'E20' = '  +0.30%  '
    'E21' = '  +0.80%  '
    'D22' = '73.34'
    'E22' = '  +0.94%  '
    'D23' = '0.998'
    'E23' = '  -0.01%  '
    'D24' = '0.535'
    'E24' = '  +0.78%  '
    'B25' = 'LEO'
    'C25' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D25' = '5.68'
    'E25' = '  -0.98%  '
    'B26' = 'PEPE'
    'C26' = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    'D26' = '0.0000121'
    'E26' = '  -0.01%  '
    'B27' = 'InternetComputer(DFINITY)'
    'C27' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D27' = '10.13'
    'E27' = '  -0.50%  '
    'B28' = 'Kaspa'
    'C28' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D28' = '0.179'
    'E28' = '  +0.38%  '
    'B29' = 'Binance-PegBSC-USD'
    'C29' = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
    'D29' = '0.997'
    'E29' = '  -0.25%  '
    'B30' = 'NEARProtocol'
    'C30' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D30' = '6.15'
    'E30' = '  -2.68%  '
    'B31' = 'Fetch.AI'
    'C31' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D31' = '1.42'
    'E31' = '  -1.56%  '
    'B32' = 'PancakeSwap'
    'C32' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D32' = '2.06'
    'E32' = '  +0.62%  '
    'B33' = 'EthereumClassic'
    'C33' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D33' = '23.67'
    'E33' = '  +0.30%  '
    'B34' = 'Aptos'
    'C34' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D34' = '7.39'
    'E34' = '  +1.14%  '
    'B35' = 'ImmutableX'
    'C35' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D35' = '1.61'
    'E35' = '  +1.52%  '
    'B36' = 'Monero'
    'C36' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D36' = '164.29'
    'E36' = '  +0.98%  '
    'B37' = 'Mantle'
    'C37' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D37' = '0.877'
    'E37' = '  -2.29%  '
    'B38' = 'Stacks'
    'C38' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'D38' = '1.90'
    'E38' = '  -0.04%  '
    'B39' = 'RenderToken'
    'C39' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D39' = '6.86'
    'E39' = '  +1.04%  '
    'B40' = 'Filecoin'
    'C40' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D40' = '4.62'
    'E40' = '  +0.07%  '
    'D41' = '2.837.63'
    'E41' = '  +2.05%  '
    'B42' = 'InjectiveProtocol'
    'C42' = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    'D42' = '27.12'
    'E42' = '  +2.21%  '
    'B43' = 'Hedera'
    'C43' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D43' = '0.0734'
    'E43' = '  -0.55%  '
    'B44' = 'EnergySwap'
    'C44' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D44' = '26.02'
    'E44' = '  -0.84%  '
    'B45' = 'OKB'
    'C45' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D45' = '42.56'
    'E45' = '  -0.26%  '
    'B46' = 'dogwifhat'
    'C46' = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    'D46' = '2.54'
    'E46' = '  +1.23%  '
    'B47' = 'VeChain'
    'C47' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D47' = '0.0301'
    'E47' = '  -2.49%  '
    'B48' = 'Bittensor'
    'C48' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'D48' = '340.07'
    'E48' = '  -1.00%  '
    'B49' = 'Arweave'
    'C49' = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
    'D49' = '34.47'
    'E49' = '  +3.92%  '
    'B50' = 'ONDO'
    'C50' = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
    'D50' = '1.07'
    'E50' = '  -0.13%  '
    'E51' = '  -1.67%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $quote + $updates[$ref]
}

Write-Host "Applied 140 cell updates"